# EPBDS-9540 Support Java Name convention on Json field name generating
# in SpreadsheetResults. Rework.
#
# The field-name placeholders used inside the "someCalcSRArray_toMap"
# spreadsheet (rows 37/38 and 56/57, columns C:H on Sheet1) referenced the
# step results in lower-case ("step1"/"step2"). They must now use the
# proper (capitalized) step names "Step1"/"Step2" so that the generated
# JSON field names follow Java naming conventions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(37, 38, 56, 57)

foreach ($r in $rows) {
    $ws.Range("C$r").Value = '_res_.$Step3["Step1"]:Integer'
    $ws.Range("D$r").Value = '_res_.$Step3["Step2"]:Integer'
    $ws.Range("E$r").Value = '_res_.$Step4["Step1"]:Integer'
    $ws.Range("F$r").Value = '_res_.$Step4["Step2"]:Integer'
    $ws.Range("G$r").Value = '_res_.$Step5["Step1"]:Integer'
    $ws.Range("H$r").Value = '_res_.$Step5["Step2"]:Integer'
}

# Keep the same look & feel as the original authoring session: move the
# active selection to where the author left off.
$ws.Range("C54").Select()
